$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Septiembre de 2020 a las 12:32"

# Countries swapped ranking: Rumania now ranks above Republica Dominicana
$ws.Range("A34").Value = "Rumania"
$ws.Range("A35").Value = "Republica Dominicana"

# Updated case numbers
$ws.Cells.Item(4, 2).Value = 6710031
$ws.Cells.Item(4, 3).Value = 1573
$ws.Cells.Item(4, 4).Value = 3975097
$ws.Cells.Item(4, 5).Value = 2536401
$ws.Cells.Item(4, 7).Value = 13
$ws.Cells.Item(4, 8).Value = 198533

$ws.Cells.Item(5, 2).Value = 4850887
$ws.Cells.Item(5, 3).Value = 5884
$ws.Cells.Item(5, 5).Value = 990996
$ws.Cells.Item(5, 7).Value = 30
$ws.Cells.Item(5, 8).Value = 79784

$ws.Cells.Item(15, 2).Value = 404648
$ws.Cells.Item(15, 3).Value = 2619
$ws.Cells.Item(15, 4).Value = 348013
$ws.Cells.Item(15, 5).Value = 33322
$ws.Cells.Item(15, 7).Value = 156
$ws.Cells.Item(15, 8).Value = 23313

$ws.Cells.Item(18, 2).Value = 339332
$ws.Cells.Item(18, 3).Value = 1812
$ws.Cells.Item(18, 4).Value = 243155
$ws.Cells.Item(18, 5).Value = 91418
$ws.Cells.Item(18, 7).Value = 26
$ws.Cells.Item(18, 8).Value = 4759

$ws.Cells.Item(26, 2).Value = 221523
$ws.Cells.Item(26, 3).Value = 3141
$ws.Cells.Item(26, 4).Value = 158405
$ws.Cells.Item(26, 5).Value = 54277
$ws.Cells.Item(26, 7).Value = 118
$ws.Cells.Item(26, 8).Value = 8841

$ws.Cells.Item(34, 2).Value = 104187
$ws.Cells.Item(34, 3).Value = 692
$ws.Cells.Item(34, 4).Value = 43244
$ws.Cells.Item(34, 5).Value = 56758
$ws.Cells.Item(34, 7).Value = 22
$ws.Cells.Item(34, 8).Value = 4185

$ws.Cells.Item(35, 2).Value = 103660
$ws.Cells.Item(35, 4).Value = 77182
$ws.Cells.Item(35, 5).Value = 24510
$ws.Cells.Item(35, 8).Value = 1968

$ws.Cells.Item(40, 2).Value = 90222
$ws.Cells.Item(40, 3).Value = 476
$ws.Cells.Item(40, 4).Value = 83928
$ws.Cells.Item(40, 5).Value = 5504
$ws.Cells.Item(40, 7).Value = 10
$ws.Cells.Item(40, 8).Value = 790

$ws.Cells.Item(66, 2).Value = 38772
$ws.Cells.Item(66, 3).Value = 56
$ws.Cells.Item(66, 4).Value = 32073
$ws.Cells.Item(66, 5).Value = 5274
$ws.Cells.Item(66, 7).Value = 5
$ws.Cells.Item(66, 8).Value = 1425

$ws.Cells.Item(76, 4).Value = 23578
$ws.Cells.Item(76, 5).Value = 2298

$ws.Cells.Item(97, 2).Value = 9946
$ws.Cells.Item(97, 3).Value = 31
$ws.Cells.Item(97, 4).Value = 9203
$ws.Cells.Item(97, 5).Value = 615

$ws.Cells.Item(179, 4).Value = 412
$ws.Cells.Item(179, 5).Value = 11

$ws.Cells.Item(183, 4).Value = 303
$ws.Cells.Item(183, 5).Value = 27
